$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 5532
$ws.Range("E2").Value = 326
$ws.Range("F2").Value = 326
$ws.Range("G2").Value = 114
$ws.Range("H2").Value = 47
$ws.Range("I2").Value = 5
$ws.Range("J2").Value = 41
$ws.Range("K2").Value = 5426
$ws.Range("L2").Value = 1451
$ws.Range("M2").Value = 3974
$ws.Range("N2").Value = 3591
$ws.Range("O2").Value = 383
$ws.Range("P2").Value = 137
$ws.Range("Q2").Value = 695
$ws.Range("R2").Value = -594
$ws.Range("S2").Value = -157
$ws.Range("T2").Value = 440
$ws.Range("U2").Value = 256
$ws.Range("V2").Value = 354
$ws.Range("W2").Value = 5.89
$ws.Range("X2").Value = 0.84
$ws.Range("Y2").Value = 0.15
$ws.Range("Z2").Value = 0.85
$ws.Range("AA2").Value = 36.52
$ws.Range("AB2").Value = 2511.5
$ws.Range("AC2").Value = 20
$ws.Range("AD2").Value = 603.42
$ws.Range("AE2").Value = 13119
$ws.Range("AF2").Value = 0.91
$ws.Range("AG2").Value = 100
$ws.Range("AH2").Value = 0.83
$ws.Range("AI2").Value = 506.29
$ws.Range("AJ2").Value = 23620751

# Row 3
$ws.Range("D3").Value = 5643
$ws.Range("E3").Value = 470
$ws.Range("F3").Value = 470
$ws.Range("G3").Value = 168
$ws.Range("H3").Value = 74
$ws.Range("I3").Value = 40
$ws.Range("J3").Value = 34
$ws.Range("K3").Value = 5270
$ws.Range("L3").Value = 1282
$ws.Range("M3").Value = 3988
$ws.Range("N3").Value = 3584
$ws.Range("O3").Value = 404
$ws.Range("P3").Value = 137
$ws.Range("Q3").Value = 616
$ws.Range("R3").Value = -320
$ws.Range("S3").Value = -28
$ws.Range("T3").Value = 337
$ws.Range("U3").Value = 279
$ws.Range("V3").Value = 356
$ws.Range("W3").Value = 8.34
$ws.Range("X3").Value = 1.31
$ws.Range("Y3").Value = 1.1
$ws.Range("Z3").Value = 1.39
$ws.Range("AA3").Value = 32.15
$ws.Range("AB3").Value = 2513.91
$ws.Range("AC3").Value = 145
$ws.Range("AD3").Value = 74.25
$ws.Range("AE3").Value = 13095
$ws.Range("AF3").Value = 0.82
$ws.Range("AG3").Value = 100
$ws.Range("AH3").Value = 0.93
$ws.Range("AI3").Value = 69.54000000000001
$ws.Range("AJ3").Value = 23620751

# Row 4
$ws.Range("D4").Value = 5453
$ws.Range("E4").Value = 160
$ws.Range("F4").Value = 160
$ws.Range("G4").Value = -56
$ws.Range("H4").Value = -42
$ws.Range("I4").Value = -68
$ws.Range("J4").Value = 26
$ws.Range("K4").Value = 5205
$ws.Range("L4").Value = 1296
$ws.Range("M4").Value = 3909
$ws.Range("N4").Value = 3483
$ws.Range("O4").Value = 426
$ws.Range("P4").Value = 137
$ws.Range("Q4").Value = 148
$ws.Range("R4").Value = -369
$ws.Range("S4").Value = -5
$ws.Range("T4").Value = 167
$ws.Range("U4").Value = -19
$ws.Range("V4").Value = 379
$ws.Range("W4").Value = 2.94
$ws.Range("X4").Value = -0.77
$ws.Range("Y4").Value = -1.92
$ws.Range("Z4").Value = -0.8
$ws.Range("AA4").Value = 33.15
$ws.Range("AB4").Value = 2445.35
$ws.Range("AC4").Value = -248
$ws.Range("AD4").Value = -45.85
$ws.Range("AE4").Value = 12727
$ws.Range("AF4").Value = 0.89
$ws.Range("AG4").Value = 100
$ws.Range("AH4").Value = 0.88
$ws.Range("AI4").Value = -40.67
$ws.Range("AJ4").Value = 23620751

# Row 5
$ws.Range("D5").Value = 5916
$ws.Range("E5").Value = 111
$ws.Range("F5").Value = 111
$ws.Range("G5").Value = 378
$ws.Range("H5").Value = 296
$ws.Range("I5").Value = 234
$ws.Range("J5").Value = 62
$ws.Range("K5").Value = 5565
$ws.Range("L5").Value = 1411
$ws.Range("M5").Value = 4154
$ws.Range("N5").Value = 3671
$ws.Range("O5").Value = 483
$ws.Range("P5").Value = 137
$ws.Range("Q5").Value = 348
$ws.Range("R5").Value = -508
$ws.Range("S5").Value = 6
$ws.Range("T5").Value = 362
$ws.Range("U5").Value = -15
$ws.Range("V5").Value = 417
$ws.Range("W5").Value = 1.87
$ws.Range("X5").Value = 5.01
$ws.Range("Y5").Value = 6.54
$ws.Range("Z5").Value = 5.5
$ws.Range("AA5").Value = 33.98
$ws.Range("AB5").Value = 2592.46
$ws.Range("AC5").Value = 855
$ws.Range("AD5").Value = 15.43
$ws.Range("AE5").Value = 13413
$ws.Range("AF5").Value = 0.98
$ws.Range("AG5").ClearContents()
$ws.Range("AH5").ClearContents()
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 23620751

# Row 6
$ws.Range("D6").Value = 5358
$ws.Range("E6").Value = -179
$ws.Range("F6").Value = -179
$ws.Range("G6").Value = -279
$ws.Range("H6").Value = -279
$ws.Range("I6").Value = -357
$ws.Range("K6").Value = 5005
$ws.Range("L6").Value = 1127
$ws.Range("M6").Value = 3878
$ws.Range("N6").Value = 3316
$ws.Range("P6").Value = 137
$ws.Range("Q6").Value = 100
$ws.Range("R6").Value = -31
$ws.Range("S6").Value = -102
$ws.Range("T6").Value = 149
$ws.Range("U6").Value = -49
$ws.Range("V6").Value = 315
$ws.Range("W6").Value = -3.34
$ws.Range("X6").Value = -5.2
$ws.Range("Y6").Value = -10.23
$ws.Range("Z6").Value = -5.27
$ws.Range("AA6").Value = 29.07
$ws.Range("AB6").Value = 2330.66
$ws.Range("AC6").Value = -1306
$ws.Range("AD6").Value = -3.91
$ws.Range("AE6").Value = 12117
$ws.Range("AF6").Value = 0.42
$ws.Range("AG6").ClearContents()
$ws.Range("AH6").ClearContents()
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 23620751

# Row 7
$ws.Range("D7").Value = 5466
$ws.Range("E7").Value = 148
$ws.Range("G7").Value = 242
$ws.Range("H7").Value = 204
$ws.Range("I7").Value = 182
$ws.Range("K7").Value = 5054
$ws.Range("L7").Value = 974
$ws.Range("M7").Value = 4080
$ws.Range("N7").Value = 3450
$ws.Range("P7").Value = 138
$ws.Range("Q7").Value = 242
$ws.Range("R7").Value = -42
$ws.Range("S7").Value = -199
$ws.Range("T7").Value = 80
$ws.Range("U7").Value = 372
$ws.Range("W7").Value = 2.7
$ws.Range("X7").Value = 3.74
$ws.Range("Y7").Value = 5.38
$ws.Range("Z7").Value = 4.07
$ws.Range("AA7").Value = 23.87
$ws.Range("AC7").Value = 665
$ws.Range("AD7").Value = 14.93
$ws.Range("AE7").Value = 12605
$ws.Range("AF7").Value = 0.79
$ws.Range("AG7").Value = 0
$ws.Range("AH7").Value = 0
$ws.Range("AI7").ClearContents()

# Row 8
$ws.Range("D8").Value = 6254
$ws.Range("E8").Value = 349
$ws.Range("G8").Value = 542
$ws.Range("H8").Value = 482
$ws.Range("I8").Value = 448
$ws.Range("K8").Value = 5632
$ws.Range("L8").Value = 1079
$ws.Range("M8").Value = 4552
$ws.Range("N8").Value = 3770
$ws.Range("P8").Value = 138
$ws.Range("Q8").Value = 370
$ws.Range("R8").Value = -150
$ws.Range("S8").Value = 30
$ws.Range("T8").Value = 300
$ws.Range("U8").Value = 360
$ws.Range("W8").Value = 5.58
$ws.Range("X8").Value = 7.7
$ws.Range("Y8").Value = 12.41
$ws.Range("Z8").Value = 9.01
$ws.Range("AA8").Value = 23.7
$ws.Range("AC8").Value = 1637
$ws.Range("AD8").Value = 6.07
$ws.Range("AE8").Value = 13775
$ws.Range("AF8").Value = 0.72
$ws.Range("AG8").Value = 0
$ws.Range("AH8").Value = 0
$ws.Range("AI8").ClearContents()

# Row 9
$ws.Range("D9").Value = 6431
$ws.Range("E9").Value = 360
$ws.Range("G9").Value = 553
$ws.Range("H9").Value = 492
$ws.Range("I9").Value = 458
$ws.Range("K9").Value = 6169
$ws.Range("L9").Value = 1129
$ws.Range("M9").Value = 5040
$ws.Range("N9").Value = 4100
$ws.Range("P9").Value = 138
$ws.Range("Q9").Value = 444
$ws.Range("R9").Value = -121
$ws.Range("S9").Value = 20
$ws.Range("T9").Value = 200
$ws.Range("U9").Value = 414
$ws.Range("W9").Value = 5.61
$ws.Range("X9").Value = 7.66
$ws.Range("Y9").Value = 11.63
$ws.Range("Z9").Value = 8.35
$ws.Range("AA9").Value = 22.4
$ws.Range("AC9").Value = 1672
$ws.Range("AD9").Value = 5.94
$ws.Range("AE9").Value = 14980
$ws.Range("AF9").Value = 0.66
$ws.Range("AG9").Value = 0
$ws.Range("AH9").Value = 0
$ws.Range("AI9").ClearContents()
